# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# Updates numeric price/profit columns (H-N) on each Leve sheet based on refreshed market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H19").Value = 5155.222
$ws.Range("I19").Value = 1299.25
$ws.Range("K19").Value = 1299.25
$ws.Range("M19").Value = -1124.25
$ws.Range("H40").Value = 3946.3635
$ws.Range("I40").Value = 2917.5
$ws.Range("J40").Value = 5181
$ws.Range("K40").Value = 2917.5
$ws.Range("L40").Value = 5181
$ws.Range("M40").Value = -2742.5
$ws.Range("N40").Value = -5531
$ws.Range("H43").Value = 20928.572
$ws.Range("J43").Value = 15800
$ws.Range("L43").Value = 15800
$ws.Range("N43").Value = -15938
$ws.Range("H64").Value = 7334.5
$ws.Range("I64").Value = 5000.3335
$ws.Range("K64").Value = 5000.3335
$ws.Range("M64").Value = -4752.3335
$ws.Range("H67").Value = 7334.5
$ws.Range("I67").Value = 5000.3335
$ws.Range("K67").Value = 5000.3335
$ws.Range("M67").Value = -4142.3335
$ws.Range("H106").Value = 4482.5
$ws.Range("I106").Value = 6497.5
$ws.Range("K106").Value = 6497.5
$ws.Range("M106").Value = -5866.5
$ws.Range("H115").Value = 7637497.5
$ws.Range("I115").Value = 9545834
$ws.Range("J115").Value = 4151.3335
$ws.Range("K115").Value = 28637502
$ws.Range("L115").Value = 12454.0005
$ws.Range("M115").Value = -28635935
$ws.Range("N115").Value = -15588.0005
$ws.Range("H132").Value = 1070
$ws.Range("I132").Value = 517.6579
$ws.Range("K132").Value = 1552.9737
$ws.Range("M132").Value = 977.0263
$ws.Range("H141").Value = 1509.625
$ws.Range("I141").Value = 1083.6
$ws.Range("J141").Value = 7900
$ws.Range("K141").Value = 3250.8
$ws.Range("L141").Value = 23700
$ws.Range("M141").Value = 1929.2
$ws.Range("N141").Value = -34060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10102853
$ws.Range("I2").Value = 12987796
$ws.Range("K2").Value = 12987796
$ws.Range("M2").Value = -12987683
$ws.Range("H32").Value = 43246.234
$ws.Range("I32").Value = 44192.43
$ws.Range("K32").Value = 44192.43
$ws.Range("M32").Value = -43905.43
$ws.Range("H45").Value = 1978.7142
$ws.Range("I45").Value = 1905.75
$ws.Range("J45").Value = 2076
$ws.Range("K45").Value = 1905.75
$ws.Range("L45").Value = 2076
$ws.Range("M45").Value = -1528.75
$ws.Range("N45").Value = -2830
$ws.Range("H74").Value = 2157.0938
$ws.Range("I74").Value = 1170.85
$ws.Range("K74").Value = 1170.85
$ws.Range("M74").Value = -296.8499999999999
$ws.Range("H77").Value = 2157.0938
$ws.Range("I77").Value = 1170.85
$ws.Range("K77").Value = 5854.25
$ws.Range("M77").Value = -1486.25
$ws.Range("H102").Value = 4431.4736
$ws.Range("I102").Value = 4431.4736
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4431.4736
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2809.4736
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 10102853
$ws.Range("I116").Value = 12987796
$ws.Range("K116").Value = 12987796
$ws.Range("M116").Value = -12985502
$ws.Range("H122").Value = 3547.6
$ws.Range("H132").Value = 5688.9375
$ws.Range("I132").Value = 4181.489
$ws.Range("J132").Value = 9259.210999999999
$ws.Range("K132").Value = 12544.467
$ws.Range("L132").Value = 27777.633
$ws.Range("M132").Value = -10014.467
$ws.Range("N132").Value = -32837.633

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10102853
$ws.Range("I3").Value = 12987796
$ws.Range("K3").Value = 12987796
$ws.Range("M3").Value = -12987682
$ws.Range("H86").Value = 54133.633
$ws.Range("I86").Value = 1495.862
$ws.Range("J86").Value = 223744.22
$ws.Range("K86").Value = 1495.862
$ws.Range("L86").Value = 223744.22
$ws.Range("M86").Value = -372.8620000000001
$ws.Range("N86").Value = -225990.22
$ws.Range("H89").Value = 54133.633
$ws.Range("I89").Value = 1495.862
$ws.Range("J89").Value = 223744.22
$ws.Range("K89").Value = 7479.31
$ws.Range("L89").Value = 1118721.1
$ws.Range("M89").Value = -1863.31
$ws.Range("N89").Value = -1129953.1
$ws.Range("H134").Value = 6686.0293
$ws.Range("I134").Value = 5055.391
$ws.Range("J134").Value = 10095.546
$ws.Range("K134").Value = 15166.173
$ws.Range("L134").Value = 30286.638
$ws.Range("M134").Value = -12631.173
$ws.Range("N134").Value = -35356.638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2482.4285
$ws.Range("I94").Value = 5037.5
$ws.Range("J94").Value = 910.0769
$ws.Range("K94").Value = 5037.5
$ws.Range("L94").Value = 910.0769
$ws.Range("M94").Value = -4586.5
$ws.Range("N94").Value = -1812.0769
$ws.Range("H105").Value = 1673
$ws.Range("I105").Value = 997
$ws.Range("K105").Value = 997
$ws.Range("M105").Value = 750
$ws.Range("H107").Value = 479.22223
$ws.Range("I107").Value = 439.125
$ws.Range("K107").Value = 439.125
$ws.Range("M107").Value = 1480.875
$ws.Range("H134").Value = 3995.7368
$ws.Range("I134").Value = 2323.125
$ws.Range("K134").Value = 6969.375
$ws.Range("M134").Value = -4434.375
$ws.Range("H141").Value = 397331.75
$ws.Range("J141").Value = 425535.56
$ws.Range("L141").Value = 425535.56
$ws.Range("N141").Value = -435895.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 921
$ws.Range("I22").Value = 500.5
$ws.Range("J22").Value = 1201.3334
$ws.Range("K22").Value = 1501.5
$ws.Range("L22").Value = 3604.0002
$ws.Range("M22").Value = -1332.5
$ws.Range("N22").Value = -3942.0002
$ws.Range("H27").Value = 921
$ws.Range("I27").Value = 500.5
$ws.Range("J27").Value = 1201.3334
$ws.Range("K27").Value = 1501.5
$ws.Range("L27").Value = 3604.0002
$ws.Range("M27").Value = -1399.5
$ws.Range("N27").Value = -3808.0002
$ws.Range("H35").Value = 1537.25
$ws.Range("I35").Value = 49.666668
$ws.Range("J35").Value = 6000
$ws.Range("K35").Value = 149.000004
$ws.Range("L35").Value = 18000
$ws.Range("M35").Value = 138.999996
$ws.Range("N35").Value = -18576
$ws.Range("H131").Value = 14499680
$ws.Range("J131").Value = 9601.75
$ws.Range("L131").Value = 28805.25
$ws.Range("N131").Value = -38885.25
$ws.Range("H137").Value = 6937.421
$ws.Range("J137").Value = 1904
$ws.Range("L137").Value = 5712
$ws.Range("N137").Value = -15912
$ws.Range("H140").Value = 1621.7
$ws.Range("J140").Value = 1717.5
$ws.Range("L140").Value = 5152.5
$ws.Range("N140").Value = -15512.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3601.7856
$ws.Range("I80").Value = 3141.5
$ws.Range("J80").Value = 3947
$ws.Range("K80").Value = 3141.5
$ws.Range("L80").Value = 3947
$ws.Range("M80").Value = -2143.5
$ws.Range("N80").Value = -5943
$ws.Range("H83").Value = 3601.7856
$ws.Range("I83").Value = 3141.5
$ws.Range("J83").Value = 3947
$ws.Range("K83").Value = 15707.5
$ws.Range("L83").Value = 19735
$ws.Range("M83").Value = -10715.5
$ws.Range("N83").Value = -29719
$ws.Range("H102").Value = 856.6667
$ws.Range("I102").Value = 635
$ws.Range("K102").Value = 635
$ws.Range("M102").Value = 987
$ws.Range("H132").Value = 4270.619
$ws.Range("J132").Value = 13503.5
$ws.Range("L132").Value = 40510.5
$ws.Range("N132").Value = -45570.5
$ws.Range("H135").Value = 60999
$ws.Range("J135").Value = 60999
$ws.Range("L135").Value = 60999
$ws.Range("N135").Value = -71139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3460.6584
$ws.Range("J22").Value = 4084.24
$ws.Range("L22").Value = 4084.24
$ws.Range("N22").Value = -4674.24
$ws.Range("H27").Value = 3460.6584
$ws.Range("J27").Value = 4084.24
$ws.Range("L27").Value = 4084.24
$ws.Range("N27").Value = -4298.24
$ws.Range("H55").Value = 205.53847
$ws.Range("I55").Value = 131.5
$ws.Range("J55").Value = 269
$ws.Range("K55").Value = 131.5
$ws.Range("L55").Value = 269
$ws.Range("M55").Value = 41.5
$ws.Range("N55").Value = -615
$ws.Range("H122").Value = 5014.952
$ws.Range("I122").Value = 4773.75
$ws.Range("K122").Value = 14321.25
$ws.Range("M122").Value = -11871.25
$ws.Range("H132").Value = 3360.96
$ws.Range("I132").Value = 2054.75
$ws.Range("J132").Value = 5683.1113
$ws.Range("K132").Value = 6164.25
$ws.Range("L132").Value = 17049.3339
$ws.Range("M132").Value = -3634.25
$ws.Range("N132").Value = -22109.3339
$ws.Range("H136").Value = 4254.9487
$ws.Range("I136").Value = 3974.9285
$ws.Range("J136").Value = 4967.727
$ws.Range("K136").Value = 11924.7855
$ws.Range("L136").Value = 14903.181
$ws.Range("M136").Value = -9374.7855
$ws.Range("N136").Value = -20003.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1026.2778
$ws.Range("I113").Value = 1196.1538
$ws.Range("K113").Value = 3588.4614
$ws.Range("M113").Value = -1418.4614
